$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'38.779.39"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "'2.104.06"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'227.62"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("D7").Value = "'62.43"
$ws.Range("E7").Value = "  +2.52%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +2.15%  "

$ws.Range("E10").Value = "  +1.08%  "

$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("D12").Value = "'15.82"
$ws.Range("E12").Value = "  +6.58%  "

$ws.Range("D13").Value = "'2.416.64"
$ws.Range("E13").Value = "  +0.70%  "

$ws.Range("D14").Value = "'22.00"
$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("D15").Value = "'0.810"
$ws.Range("E15").Value = "  +3.10%  "

$ws.Range("E16").Value = "  +1.54%  "

$ws.Range("D17").Value = "'2.102.82"
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").Value = "'38.796.46"
$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("D19").Value = "'6.11"
$ws.Range("E19").Value = "  +0.78%  "

$ws.Range("D20").Value = "'71.51"
$ws.Range("E20").Value = "  +0.32%  "

$ws.Range("D21").Value = "'0.0₃0848"
$ws.Range("E21").Value = "  +1.67%  "

$ws.Range("D22").Value = "'228.15"
$ws.Range("E22").Value = "  +1.13%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").Value = "'2.35"
$ws.Range("E24").Value = "  -1.77%  "

$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("D26").Value = "'9.65"
$ws.Range("E26").Value = "  +2.33%  "

$ws.Range("D27").Value = "'172.39"
$ws.Range("E27").Value = "  +1.39%  "

$ws.Range("D28").Value = "'0.137"
$ws.Range("E28").Value = "  +1.61%  "

$ws.Range("E29").Value = "  +3.07%  "

$ws.Range("E30").Value = "  +1.49%  "

$ws.Range("D31").Value = "'2.57"
$ws.Range("E31").Value = "  +9.90%  "

$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("E33").Value = "  +1.13%  "

$ws.Range("D34").Value = "'7.17"
$ws.Range("E34").Value = "  +11.19%  "

$ws.Range("E35").Value = "  -1.14%  "

$ws.Range("E36").Value = "  +2.03%  "

$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").Value = "'3.53"
$ws.Range("E38").Value = "  -0.99%  "

$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("D40").Value = "'18.07"
$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("D41").Value = "'102.36"
$ws.Range("E41").Value = "  +2.31%  "

$ws.Range("E42").Value = "  +3.44%  "

$ws.Range("D43").Value = "'1.525.79"
$ws.Range("E43").Value = "  -0.94%  "

$ws.Range("D44").Value = "'1.20"
$ws.Range("E44").Value = "  +7.73%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'7.84"
$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "'2.81"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "'0.0915"
$ws.Range("E47").Value = "  -1.57%  "

$ws.Range("D48").Value = "'4.21"
$ws.Range("E48").Value = "  +1.81%  "

$ws.Range("E49").Value = "  +4.09%  "

$ws.Range("E50").Value = "  -0.82%  "

$ws.Range("D51").Value = "'2.303.89"
$ws.Range("E51").Value = "  +0.81%  "
